$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the Recall / Precision formulas, which had been swapped ---
# Recall  = TP / (TP + FN)   -> H4 / (H4 + H3)
# Precision = TP / (TP + FP) -> H4 / (H4 + G4)
$ws.Range("D10").Formula = "=H4/(H4+H3)"
$ws.Range("D11").Formula = "=H4/(H4+G4)"

# --- Fix the matching helper-text annotations in column E (also swapped) ---
$ws.Range("E10").Value = " (TP) / (TP + FP) - True positive rate / Sensitivity / Prob. Of detection"
$ws.Range("E11").Value = " (TP) / (TP + FN) - Fraction of +ive predictions as correct"

# Re-apply the original cell formatting to E10:E11 (writing .Value resets it)
$ws.Range("E13").Copy()
$ws.Range("E10:E11").PasteSpecial(-4122)  # xlPasteFormats

# --- Remove leftover scratch / helper cells that are no longer needed ---
$ws.Range("K8").Clear()
$ws.Range("H17").Clear()
$ws.Range("G17").Value = ""
$ws.Range("G18").Value = ""
$ws.Range("H18").Value = ""
$ws.Rows(19).Delete()

# --- Update the saved view/selection ---
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("F16").Select()
